$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Schedule")
$ws.Range("A3").Value = 46053.27083333334
$ws.Range("B3").Value = 46053.77083333334
$ws.Range("E3").Value = 845.3962432500002
$ws.Range("F3").Value = 18.63748331679895

$ws = $wb.Worksheets.Item("Detailed")
$ws.Range("B37").Value = 101.25
$ws.Range("B38").Value = 108.89
$ws.Range("B39").Value = 123.17919
$ws.Range("B40").Value = 283.96
$ws.Range("C40").Value = "historical"
$ws.Range("B41").Value = 299.99
$ws.Range("C41").Value = "historical"
$ws.Range("B42").Value = 299.99
$ws.Range("C42").Value = "historical"
$ws.Range("B43").Value = 299.98
$ws.Range("C43").Value = "historical"
$ws.Range("B44").Value = 240.89
$ws.Range("C44").Value = "historical"
$ws.Range("B45").Value = 147.89
$ws.Range("C45").Value = "historical"
$ws.Range("B46").Value = 138.2218
$ws.Range("C46").Value = "historical"
$ws.Range("B47").Value = 109.51175
$ws.Range("C47").Value = "historical"
$ws.Range("B48").Value = 112.49854
$ws.Range("C48").Value = "historical"
$ws.Range("B49").Value = 105.79
$ws.Range("B50").Value = 105.79
$ws.Range("B51").Value = 105.79
$ws.Range("B52").Value = 108.89
$ws.Range("B53").Value = 105.12588
$ws.Range("B54").Value = 105
$ws.Range("B55").Value = 105.78998
$ws.Range("B56").Value = 108.89
$ws.Range("B57").Value = 108.89
$ws.Range("B58").Value = 108.89
$ws.Range("B59").Value = 108.89
$ws.Range("B60").Value = 108.89
$ws.Range("B61").Value = 108.89
$ws.Range("B62").Value = 108.89
$ws.Range("E62").Value = "OFF"
$ws.Range("B63").Value = 84.79000000000001
$ws.Range("B64").Value = 58.18999
$ws.Range("B65").Value = 39.93361
$ws.Range("B67").Value = 35.53145
$ws.Range("B68").Value = 0.51
$ws.Range("B71").Value = 22.07
$ws.Range("B72").Value = 13.11547
$ws.Range("B73").Value = 10.40562
$ws.Range("B74").Value = 21.42291
$ws.Range("B75").Value = 22.07
$ws.Range("B76").Value = 9.63969
$ws.Range("B77").Value = 9.48333
$ws.Range("B78").Value = 10.34407
$ws.Range("B79").Value = 36.0601
$ws.Range("B80").Value = 56.98
$ws.Range("B81").Value = 57.57846
$ws.Range("B82").Value = 50.27428
$ws.Range("B83").Value = 55.78964
$ws.Range("B84").Value = 57.35
$ws.Range("B85").Value = 73.45435000000001
$ws.Range("B86").Value = 105
$ws.Range("E86").Value = "ON"
$ws.Range("B87").Value = 105
$ws.Range("B88").Value = 147.52
$ws.Range("B89").Value = 147.52
$ws.Range("B90").Value = 132.6472
$ws.Range("B91").Value = 107.54808
$ws.Range("B93").Value = 105.79
$ws.Range("B94").Value = 100.3
$ws.Range("B95").Value = 105
$ws.Range("B96").Value = 105.79
$ws.Range("B97").Value = 97.42995000000001

